# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets share the identical data set, so the same row -> new value
# mapping applies to both of them.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 56
    4  = 103
    5  = 139
    6  = 1348
    7  = 1580
    8  = 347
    9  = 448
    11 = 179
    15 = 286
    16 = 323
    18 = 1772
    22 = 699
    24 = 348
    25 = 4288
    27 = 294
    28 = 1128
    29 = 496
    31 = 645
    33 = 327
    34 = 51
    35 = 166
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
